# Generate Report for Handback
#
# The 0141382f-82d5-4753-a969-633695d4b7d2.md file has been handed back
# (in sync with en-US). Its row moves to the top of each status table and
# picks up target/handback file info + timestamps; the two rows that used
# to precede it (ffffc185753e...md and ffffff586ce940...md) shift down by
# one position, keeping their own data intact.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "0141382f-82d5-4753-a969-633695d4b7d2.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-21 10:59:53"

$ws.Range("A3").Value = "ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-03-21 10:58:35"

$ws.Range("A4").Value = "ffffff586ce940-6904-4f5b-8777-5df999880bed.md"
$ws.Range("B4").Value = "Handed back: in sync with en-US"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "2016-03-21 10:58:35"

# Hyperlinks.Delete() removes every hyperlink on the sheet, so rebuild all
# three (the underlying target URLs are unchanged - only which display
# text/cell they are attached to changes).
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/583e461fca958a37d4637db29409c481b135fa30/e2e/ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md", "", "", "0141382f-82d5-4753-a969-633695d4b7d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/583e461fca958a37d4637db29409c481b135fa30/e2e/ffffff586ce940-6904-4f5b-8777-5df999880bed.md", "", "", "ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/583e461fca958a37d4637db29409c481b135fa30/e2e/0141382f-82d5-4753-a969-633695d4b7d2.md", "", "", "ffffff586ce940-6904-4f5b-8777-5df999880bed.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "0141382f-82d5-4753-a969-633695d4b7d2.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "0141382f-82d5-4753-a969-633695d4b7d2.a270409dde14787869912a1b2fdc7f361f0f03e3.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-21 10:59:49"
$ws.Range("F2").Value = "0141382f-82d5-4753-a969-633695d4b7d2.md"
$ws.Range("G2").Value = "0141382f-82d5-4753-a969-633695d4b7d2.a270409dde14787869912a1b2fdc7f361f0f03e3.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-21 11:00:20"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-21 10:58:32"
$ws.Range("F3").Value = "50527308-acf1-477c-ac29-3589133d0d67.md"
$ws.Range("G3").Value = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-21 10:58:52"
$ws.Range("J3").Value = "Include"

$ws.Range("A4").Value = "ffffff586ce940-6904-4f5b-8777-5df999880bed.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-21 10:58:32"
$ws.Range("F4").Value = "50527308-acf1-477c-ac29-3589133d0d67.md"
$ws.Range("G4").Value = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf"
$ws.Range("H4").Value = "2016-03-21 10:58:52"
$ws.Range("J4").Value = "Include"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/583e461fca958a37d4637db29409c481b135fa30/e2e/ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md", "", "", "0141382f-82d5-4753-a969-633695d4b7d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e584c33e8c2f998e3c3789aa24b77cf96e36fe55/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf", "", "", "0141382f-82d5-4753-a969-633695d4b7d2.a270409dde14787869912a1b2fdc7f361f0f03e3.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f784834afbb57b9caf2b791e80af92571e5d33f0/e2e/50527308-acf1-477c-ac29-3589133d0d67.md", "", "", "0141382f-82d5-4753-a969-633695d4b7d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7715dcb2d3a4f77b51005004b5c37d473940b16f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf", "", "", "0141382f-82d5-4753-a969-633695d4b7d2.a270409dde14787869912a1b2fdc7f361f0f03e3.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/583e461fca958a37d4637db29409c481b135fa30/e2e/ffffff586ce940-6904-4f5b-8777-5df999880bed.md", "", "", "ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e584c33e8c2f998e3c3789aa24b77cf96e36fe55/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf", "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f784834afbb57b9caf2b791e80af92571e5d33f0/e2e/50527308-acf1-477c-ac29-3589133d0d67.md", "", "", "50527308-acf1-477c-ac29-3589133d0d67.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7715dcb2d3a4f77b51005004b5c37d473940b16f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf", "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/583e461fca958a37d4637db29409c481b135fa30/e2e/0141382f-82d5-4753-a969-633695d4b7d2.md", "", "", "ffffff586ce940-6904-4f5b-8777-5df999880bed.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e584c33e8c2f998e3c3789aa24b77cf96e36fe55/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf", "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f784834afbb57b9caf2b791e80af92571e5d33f0/e2e/50527308-acf1-477c-ac29-3589133d0d67.md", "", "", "50527308-acf1-477c-ac29-3589133d0d67.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7715dcb2d3a4f77b51005004b5c37d473940b16f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf", "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "0141382f-82d5-4753-a969-633695d4b7d2.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "0141382f-82d5-4753-a969-633695d4b7d2.a270409dde14787869912a1b2fdc7f361f0f03e3.de-de.xlf"
$ws.Range("E2").Value = "2016-03-21 10:59:53"
$ws.Range("F2").Value = "0141382f-82d5-4753-a969-633695d4b7d2.md"
$ws.Range("G2").Value = "0141382f-82d5-4753-a969-633695d4b7d2.a270409dde14787869912a1b2fdc7f361f0f03e3.de-de.xlf"
$ws.Range("H2").Value = "2016-03-21 11:00:29"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf"
$ws.Range("E3").Value = "2016-03-21 10:58:35"
$ws.Range("F3").Value = "50527308-acf1-477c-ac29-3589133d0d67.md"
$ws.Range("G3").Value = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf"
$ws.Range("H3").Value = "2016-03-21 10:58:58"
$ws.Range("J3").Value = "Include"

$ws.Range("A4").Value = "ffffff586ce940-6904-4f5b-8777-5df999880bed.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf"
$ws.Range("E4").Value = "2016-03-21 10:58:35"
$ws.Range("F4").Value = "50527308-acf1-477c-ac29-3589133d0d67.md"
$ws.Range("G4").Value = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf"
$ws.Range("H4").Value = "2016-03-21 10:58:58"
$ws.Range("J4").Value = "Include"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/583e461fca958a37d4637db29409c481b135fa30/e2e/ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md", "", "", "0141382f-82d5-4753-a969-633695d4b7d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/28fac1fbaa18e9bee999b896c6a1cafbf4b5673b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf", "", "", "0141382f-82d5-4753-a969-633695d4b7d2.a270409dde14787869912a1b2fdc7f361f0f03e3.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8a1ba931542e69ed3338756e8c450a7718fba944/e2e/50527308-acf1-477c-ac29-3589133d0d67.md", "", "", "0141382f-82d5-4753-a969-633695d4b7d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0abb1900b2742a02ec072a2c12252b0f6bb416b9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf", "", "", "0141382f-82d5-4753-a969-633695d4b7d2.a270409dde14787869912a1b2fdc7f361f0f03e3.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/583e461fca958a37d4637db29409c481b135fa30/e2e/ffffff586ce940-6904-4f5b-8777-5df999880bed.md", "", "", "ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/28fac1fbaa18e9bee999b896c6a1cafbf4b5673b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf", "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8a1ba931542e69ed3338756e8c450a7718fba944/e2e/50527308-acf1-477c-ac29-3589133d0d67.md", "", "", "50527308-acf1-477c-ac29-3589133d0d67.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0abb1900b2742a02ec072a2c12252b0f6bb416b9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf", "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/583e461fca958a37d4637db29409c481b135fa30/e2e/0141382f-82d5-4753-a969-633695d4b7d2.md", "", "", "ffffff586ce940-6904-4f5b-8777-5df999880bed.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/28fac1fbaa18e9bee999b896c6a1cafbf4b5673b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf", "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8a1ba931542e69ed3338756e8c450a7718fba944/e2e/50527308-acf1-477c-ac29-3589133d0d67.md", "", "", "50527308-acf1-477c-ac29-3589133d0d67.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0abb1900b2742a02ec072a2c12252b0f6bb416b9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf", "", "", "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf") | Out-Null
